# Insert a new data row for "Ají" / variety "Inferno" (Región de Arica y
# Parinacota) right after the current row 197, which pushes every
# subsequent row (old 198..234) down by one (new 199..235).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 198..234 down to 199..235 by inserting a new blank row at 198.
$ws.Rows.Item(198).EntireRow.Insert()

# Populate the newly inserted row 198 with the new record's data.
$ws.Cells.Item(198, 1).Value  = 11
$ws.Cells.Item(198, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(198, 3).Value  = "Bíobío"
$ws.Cells.Item(198, 4).Value  = 45204
$ws.Cells.Item(198, 5).Value  = 8
$ws.Cells.Item(198, 6).Value  = 100112021
$ws.Cells.Item(198, 7).Value  = "Ají"
$ws.Cells.Item(198, 8).Value  = "Inferno"
$ws.Cells.Item(198, 9).Value  = "Primera"
$ws.Cells.Item(198, 10).Value = 50
$ws.Cells.Item(198, 11).Value = 36000
$ws.Cells.Item(198, 12).Value = 37000
$ws.Cells.Item(198, 13).Value = 36600
$ws.Cells.Item(198, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(198, 16).Value = 3660
$ws.Cells.Item(198, 17).Value = 10
$ws.Cells.Item(198, 18).Value = "Hortaliza"
